$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added for "Ajo" (garlic) at Macroferia
# Regional de Talca. It belongs chronologically right after the existing
# row 157, so insert a fresh row at 158 and push the old rows 158-179
# down to 159-180 (the old row 179 ends up as the new row 180).
$ws.Rows.Item(158).Insert()

# Populate the newly inserted row 158 with the new weekly record. Its
# values mirror the (now shifted) row below it, except for the new date.
$ws.Cells.Item(158, 1).Value = 5
$ws.Cells.Item(158, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(158, 3).Value = "Maule"
$ws.Cells.Item(158, 4).Value = 44491
$ws.Cells.Item(158, 5).Value = 7
$ws.Cells.Item(158, 6).Value = 100112003
$ws.Cells.Item(158, 7).Value = "Ajo"
$ws.Cells.Item(158, 8).Value = "Chino"
$ws.Cells.Item(158, 9).Value = "Primera"
$ws.Cells.Item(158, 10).Value = 200
$ws.Cells.Item(158, 11).Value = 14000
$ws.Cells.Item(158, 12).Value = 14000
$ws.Cells.Item(158, 13).Value = 14000
$ws.Cells.Item(158, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(158, 15).Value = "China"
$ws.Cells.Item(158, 16).Value = 1400
$ws.Cells.Item(158, 17).Value = 10
$ws.Cells.Item(158, 18).Value = "Hortaliza"

# Make sure the new row's date cell keeps the same date-time number
# format as the rest of column D.
$ws.Cells.Item(158, 4).NumberFormat = $ws.Cells.Item(159, 4).NumberFormat
